$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the "_old" / "_new" column header suffixes to the respective
#    input-file-name suffixes ("_FV2304" / "_FV2310").
# ---------------------------------------------------------------------------
$oldHeaders = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

$newHeaders = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}

# Column 11 ("diff") keeps its name.

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i]
}

# ---------------------------------------------------------------------------
# 2. Turn the data range into a real Excel Table ("Table1") so the header
#    row gets filter buttons and the data is addressable as a ListObject.
# ---------------------------------------------------------------------------
$rng = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# ---------------------------------------------------------------------------
# 3. Freeze the header row (split below row 1, freeze top-left cell A2).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
